$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C16").Value = "j'ai fait une modif!"

$ws.Range("C17").Select()
